$wb = $excel.ActiveWorkbook

# ==== Sheet: Overview ====
$ws = $wb.Worksheets.Item("Overview")

# -- set cell values --
$ws.Range("A1").Value = "File Name"
$ws.Range("B1").Value = "zh-cn"
$ws.Range("C1").Value = "de-de"
$ws.Range("A2").Value = "ffff28ca022c-c202-4afb-b2de-7d381b8a1aa4.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("A3").Value = "ffffff0bfef447-7142-48ff-85c4-1c1be77f12d9.md"
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("A4").Value = "24745f76-4497-4736-9139-e34de63a432a.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("A5").Value = ".localization-config"
$ws.Range("B5").Value = "Not to be localized"
$ws.Range("C5").Value = "Not to be localized"

# -- rebuild hyperlinks (position -> URL mapping is fixed; only display text changes) --
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/7fc304006046ab89a02b83a319f28722603a74f3/e2e/24745f76-4497-4736-9139-e34de63a432a.md", "", "", "ffff28ca022c-c202-4afb-b2de-7d381b8a1aa4.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/7fc304006046ab89a02b83a319f28722603a74f3/e2e/ffff28ca022c-c202-4afb-b2de-7d381b8a1aa4.md", "", "", "ffffff0bfef447-7142-48ff-85c4-1c1be77f12d9.md")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/7fc304006046ab89a02b83a319f28722603a74f3/e2e/ffffff0bfef447-7142-48ff-85c4-1c1be77f12d9.md", "", "", "24745f76-4497-4736-9139-e34de63a432a.md")
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/7fc304006046ab89a02b83a319f28722603a74f3/.localization-config", "", "", ".localization-config")

# ==== Sheet: zh-cn ====
$ws = $wb.Worksheets.Item("zh-cn")

# -- set cell values --
$ws.Range("A1").Value = "Source File Name"
$ws.Range("B1").Value = "Status"
$ws.Range("C1").Value = "Latest Handoff File"
$ws.Range("D1").Value = "Latest Handoff Datetime"
$ws.Range("E1").Value = "Latest Target File"
$ws.Range("F1").Value = "Latest Handback File"
$ws.Range("G1").Value = "Latest Handback DateTime"
$ws.Range("H1").Value = "Handoff Reason"
$ws.Range("I1").Value = "Dependency From"
$ws.Range("A2").Value = "ffff28ca022c-c202-4afb-b2de-7d381b8a1aa4.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.d76442c13d83cd579ed92490c6b7780c957ec87e.zh-cn.xlf"
$ws.Range("D2").Value = "2016-01-28 09:37:23"
$ws.Range("E2").Value = "4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.md"
$ws.Range("F2").Value = "4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.d76442c13d83cd579ed92490c6b7780c957ec87e.zh-cn.xlf"
$ws.Range("G2").Value = "2016-01-28 09:38:09"
$ws.Range("H2").Value = "Include"
$ws.Range("A3").Value = "ffffff0bfef447-7142-48ff-85c4-1c1be77f12d9.md"
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("C3").Value = "4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.d76442c13d83cd579ed92490c6b7780c957ec87e.zh-cn.xlf"
$ws.Range("D3").Value = "2016-01-28 09:37:23"
$ws.Range("E3").Value = "4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.md"
$ws.Range("F3").Value = "4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.d76442c13d83cd579ed92490c6b7780c957ec87e.zh-cn.xlf"
$ws.Range("G3").Value = "2016-01-28 09:38:09"
$ws.Range("H3").Value = "Include"
$ws.Range("A4").Value = "24745f76-4497-4736-9139-e34de63a432a.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "24745f76-4497-4736-9139-e34de63a432a.cef8b5635807256dfb783ebb223c768ad826ab81.zh-cn.xlf"
$ws.Range("D4").Value = "2016-01-28 09:41:58"
$ws.Range("E4").Value = "24745f76-4497-4736-9139-e34de63a432a.md"
$ws.Range("F4").Value = "24745f76-4497-4736-9139-e34de63a432a.cef8b5635807256dfb783ebb223c768ad826ab81.zh-cn.xlf"
$ws.Range("G4").Value = "2016-01-28 09:40:47"
$ws.Range("H4").Value = "Include"
$ws.Range("A5").Value = ".localization-config"
$ws.Range("B5").Value = "Not to be localized"
$ws.Range("D5").Value = "0001-01-01 00:00:00"
$ws.Range("G5").Value = "0001-01-01 00:00:00"
$ws.Range("H5").Value = "Ignored"

# -- rebuild hyperlinks (position -> URL mapping is fixed; only display text changes) --
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/7fc304006046ab89a02b83a319f28722603a74f3/e2e/24745f76-4497-4736-9139-e34de63a432a.md", "", "", "ffff28ca022c-c202-4afb-b2de-7d381b8a1aa4.md")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bc83164bdf673c151a725cd9f811b946100045e5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/24745f76-4497-4736-9139-e34de63a432a.cef8b5635807256dfb783ebb223c768ad826ab81.zh-cn.xlf", "", "", "4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.d76442c13d83cd579ed92490c6b7780c957ec87e.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/98ee5df8f96bc356eef088fe3020bea33f8b26e2/e2e/24745f76-4497-4736-9139-e34de63a432a.md", "", "", "4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.md")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/09b97dcbfa22370f5a2627e65411d7523926a594/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/24745f76-4497-4736-9139-e34de63a432a.cef8b5635807256dfb783ebb223c768ad826ab81.zh-cn.xlf", "", "", "4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.d76442c13d83cd579ed92490c6b7780c957ec87e.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/7fc304006046ab89a02b83a319f28722603a74f3/e2e/ffff28ca022c-c202-4afb-b2de-7d381b8a1aa4.md", "", "", "ffffff0bfef447-7142-48ff-85c4-1c1be77f12d9.md")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/daa45f64cde85a6f199e279e900d5e46f565234b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.d76442c13d83cd579ed92490c6b7780c957ec87e.zh-cn.xlf", "", "", "4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.d76442c13d83cd579ed92490c6b7780c957ec87e.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/81ae5652ca9655c165586ddfba6000886c4a3da2/e2e/4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.md", "", "", "4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.md")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3b723ca59f7bfc22e775d035e3b65baacd1a6202/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.d76442c13d83cd579ed92490c6b7780c957ec87e.zh-cn.xlf", "", "", "4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.d76442c13d83cd579ed92490c6b7780c957ec87e.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/7fc304006046ab89a02b83a319f28722603a74f3/e2e/ffffff0bfef447-7142-48ff-85c4-1c1be77f12d9.md", "", "", "24745f76-4497-4736-9139-e34de63a432a.md")
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/daa45f64cde85a6f199e279e900d5e46f565234b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.d76442c13d83cd579ed92490c6b7780c957ec87e.zh-cn.xlf", "", "", "24745f76-4497-4736-9139-e34de63a432a.cef8b5635807256dfb783ebb223c768ad826ab81.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/81ae5652ca9655c165586ddfba6000886c4a3da2/e2e/4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.md", "", "", "24745f76-4497-4736-9139-e34de63a432a.md")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3b723ca59f7bfc22e775d035e3b65baacd1a6202/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.d76442c13d83cd579ed92490c6b7780c957ec87e.zh-cn.xlf", "", "", "24745f76-4497-4736-9139-e34de63a432a.cef8b5635807256dfb783ebb223c768ad826ab81.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/7fc304006046ab89a02b83a319f28722603a74f3/.localization-config", "", "", ".localization-config")

# ==== Sheet: de-de ====
$ws = $wb.Worksheets.Item("de-de")

# -- set cell values --
$ws.Range("A1").Value = "Source File Name"
$ws.Range("B1").Value = "Status"
$ws.Range("C1").Value = "Latest Handoff File"
$ws.Range("D1").Value = "Latest Handoff Datetime"
$ws.Range("E1").Value = "Latest Target File"
$ws.Range("F1").Value = "Latest Handback File"
$ws.Range("G1").Value = "Latest Handback DateTime"
$ws.Range("H1").Value = "Handoff Reason"
$ws.Range("I1").Value = "Dependency From"
$ws.Range("A2").Value = "ffff28ca022c-c202-4afb-b2de-7d381b8a1aa4.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.d76442c13d83cd579ed92490c6b7780c957ec87e.de-de.xlf"
$ws.Range("D2").Value = "2016-01-28 09:37:36"
$ws.Range("E2").Value = "4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.md"
$ws.Range("F2").Value = "4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.d76442c13d83cd579ed92490c6b7780c957ec87e.de-de.xlf"
$ws.Range("G2").Value = "2016-01-28 09:38:30"
$ws.Range("H2").Value = "Include"
$ws.Range("A3").Value = "ffffff0bfef447-7142-48ff-85c4-1c1be77f12d9.md"
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("C3").Value = "4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.d76442c13d83cd579ed92490c6b7780c957ec87e.de-de.xlf"
$ws.Range("D3").Value = "2016-01-28 09:37:36"
$ws.Range("E3").Value = "4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.md"
$ws.Range("F3").Value = "4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.d76442c13d83cd579ed92490c6b7780c957ec87e.de-de.xlf"
$ws.Range("G3").Value = "2016-01-28 09:38:30"
$ws.Range("H3").Value = "Include"
$ws.Range("A4").Value = "24745f76-4497-4736-9139-e34de63a432a.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "24745f76-4497-4736-9139-e34de63a432a.cef8b5635807256dfb783ebb223c768ad826ab81.de-de.xlf"
$ws.Range("D4").Value = "2016-01-28 09:42:11"
$ws.Range("E4").Value = "24745f76-4497-4736-9139-e34de63a432a.md"
$ws.Range("F4").Value = "24745f76-4497-4736-9139-e34de63a432a.cef8b5635807256dfb783ebb223c768ad826ab81.de-de.xlf"
$ws.Range("G4").Value = "2016-01-28 09:41:10"
$ws.Range("H4").Value = "Include"
$ws.Range("A5").Value = ".localization-config"
$ws.Range("B5").Value = "Not to be localized"
$ws.Range("D5").Value = "0001-01-01 00:00:00"
$ws.Range("G5").Value = "0001-01-01 00:00:00"
$ws.Range("H5").Value = "Ignored"

# -- rebuild hyperlinks (position -> URL mapping is fixed; only display text changes) --
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/7fc304006046ab89a02b83a319f28722603a74f3/e2e/24745f76-4497-4736-9139-e34de63a432a.md", "", "", "ffff28ca022c-c202-4afb-b2de-7d381b8a1aa4.md")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/366a28146f7be6c1f8a0054c35230253f5fee61a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/24745f76-4497-4736-9139-e34de63a432a.cef8b5635807256dfb783ebb223c768ad826ab81.de-de.xlf", "", "", "4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.d76442c13d83cd579ed92490c6b7780c957ec87e.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/544360d5e7898b8c92a80992a171e2f4ed93f747/e2e/24745f76-4497-4736-9139-e34de63a432a.md", "", "", "4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.md")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/1b085ce936aa1e4d99186c00855bd1dd2fa6cd51/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/24745f76-4497-4736-9139-e34de63a432a.cef8b5635807256dfb783ebb223c768ad826ab81.de-de.xlf", "", "", "4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.d76442c13d83cd579ed92490c6b7780c957ec87e.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/7fc304006046ab89a02b83a319f28722603a74f3/e2e/ffff28ca022c-c202-4afb-b2de-7d381b8a1aa4.md", "", "", "ffffff0bfef447-7142-48ff-85c4-1c1be77f12d9.md")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/426b6feca0c0f35facdbef193cd977f5b1d3718b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.d76442c13d83cd579ed92490c6b7780c957ec87e.de-de.xlf", "", "", "4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.d76442c13d83cd579ed92490c6b7780c957ec87e.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/bb218d3fe5ef63f6fff279bbb721173f0ab2e064/e2e/4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.md", "", "", "4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.md")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/498d7b89654cdebb4dd5705bf89a70a4b8038c89/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.d76442c13d83cd579ed92490c6b7780c957ec87e.de-de.xlf", "", "", "4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.d76442c13d83cd579ed92490c6b7780c957ec87e.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/7fc304006046ab89a02b83a319f28722603a74f3/e2e/ffffff0bfef447-7142-48ff-85c4-1c1be77f12d9.md", "", "", "24745f76-4497-4736-9139-e34de63a432a.md")
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/426b6feca0c0f35facdbef193cd977f5b1d3718b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.d76442c13d83cd579ed92490c6b7780c957ec87e.de-de.xlf", "", "", "24745f76-4497-4736-9139-e34de63a432a.cef8b5635807256dfb783ebb223c768ad826ab81.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/bb218d3fe5ef63f6fff279bbb721173f0ab2e064/e2e/4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.md", "", "", "24745f76-4497-4736-9139-e34de63a432a.md")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/498d7b89654cdebb4dd5705bf89a70a4b8038c89/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.d76442c13d83cd579ed92490c6b7780c957ec87e.de-de.xlf", "", "", "24745f76-4497-4736-9139-e34de63a432a.cef8b5635807256dfb783ebb223c768ad826ab81.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/7fc304006046ab89a02b83a319f28722603a74f3/.localization-config", "", "", ".localization-config")
